# Update the Adam15-Itgav LR-pair sheet with the refreshed TPM-derived
# numbers ("update scripts wuth new tpm").
#
# The upstream NATMI scoring script recomputed average/total expression
# values, per-cluster derived specificity, and edge weights after the TPM
# values were regenerated. Concretely, only the "ECs" cluster's average
# ligand (ECs sending Adam15) and average receptor (ECs receiving Itgav)
# expression values changed; that ripples into every derived column
# (total expression, derived specificity, edge weight, edge specificity)
# across all 9 data rows because those are all normalised against the
# per-cluster totals. The values below are exactly the refreshed numbers
# produced by that rerun.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: G,H,I,J = ligand avg/total/avgSpec/totalSpec
#          M,N,O,P = receptor avg/total/avgSpec/totalSpec
#          Q,R,S,T = edge avg weight/total weight/avg spec/total spec
$rowValues = @{
    2  = @{ G=31.749936;        H=95.249808;        I=0.5302851438878331; J=0.5302851438878331;
            M=3.759736666666667; N=11.27921;         O=0.0683751702595819;  P=0.06837517025958188;
            Q=119.37139854352;   R=1074.34258689168;  S=0.03625833699945747; T=0.03625833699945746 }
    3  = @{ G=31.749936;        H=95.249808;        I=0.5302851438878331; J=0.5302851438878331;
                                                      O=0.6514180024294648;  P=0.6514180024294647;
            Q=1137.264853472656; R=10235.3836812539;  S=0.3454372891494336;  T=0.3454372891494335 }
    4  = @{ G=31.749936;        H=95.249808;        I=0.5302851438878331; J=0.5302851438878331;
                                                      O=0.2802068273109533;  P=0.2802068273109533;
            Q=489.19338307408;   R=4402.74044766672;  S=0.1485895177389421;  T=0.1485895177389421 }
    5  = @{                                          I=0.3451699599880819; J=0.3451699599880819;
            M=3.759736666666667; N=11.27921;         O=0.0683751702595819;  P=0.06837517025958188;
            Q=77.70050006851332; R=699.3045006166199; S=0.02360105478267817; T=0.02360105478267817 }
    6  = @{                                          I=0.3451699599880819; J=0.3451699599880819;
                                                      O=0.6514180024294648;  P=0.6514180024294647;
                                                      S=0.2248499258340946;  T=0.2248499258340945 }
    7  = @{                                          I=0.3451699599880819; J=0.3451699599880819;
                                                      O=0.2802068273109533;  P=0.2802068273109533;
                                                      S=0.09671897937130912; T=0.0967189793713091 }
    8  = @{                                          I=0.1245448961240849; J=0.1245448961240849;
            M=3.759736666666667; N=11.27921;         O=0.0683751702595819;  P=0.06837517025958188;
            Q=28.03604551843555; R=252.32440966592;   S=0.008515778477446246; T=0.008515778477446244 }
    9  = @{                                          I=0.1245448961240849; J=0.1245448961240849;
                                                      O=0.6514180024294648;  P=0.6514180024294647;
                                                      S=0.08113078744593658; T=0.08113078744593655 }
    10 = @{                                          I=0.1245448961240849; J=0.1245448961240849;
                                                      O=0.2802068273109533;  P=0.2802068273109533;
                                                      S=0.03489833020070208; T=0.03489833020070206 }
}

$colIndex = @{ G=7; H=8; I=9; J=10; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }

foreach ($r in $rowValues.Keys) {
    $cols = $rowValues[$r]
    foreach ($colName in $cols.Keys) {
        $c = $colIndex[$colName]
        $ws.Cells.Item($r, $c).Value2 = $cols[$colName]
    }
}
